$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164479374885559
$ws.Range("B1").Value = 2.420864820480347
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.377026081085205
$ws.Range("E1").Value = 1.235227346420288
